$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1 ("I0") and J1 ("IF") — reuse H1's style (bold, bordered,
# centered header format) by copying its formatting, then set the text.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for I2:J40
$data = @{
    2  = @(3, 4)
    3  = @(4, 4)
    4  = @(3, 4)
    5  = @(7, 7)
    6  = @(7, 7)
    7  = @(1, 2)
    8  = @(5, 6)
    9  = @(7, 8)
    10 = @(9, 9)
    11 = @(6, 6)
    12 = @(9, 9)
    13 = @(8, 8)
    14 = @(7, 7)
    15 = @(8, 8)
    16 = @(7, 8)
    17 = @(7, 7)
    18 = @(8, 8)
    19 = @(7, 7)
    20 = @(10, 10)
    21 = @(5, 6)
    22 = @(8, 8)
    23 = @(7, 7)
    24 = @(5, 6)
    25 = @(7, 7)
    26 = @(8, 8)
    27 = @(9, 9)
    28 = @(5, 6)
    29 = @(7, 7)
    30 = @(4, 5)
    31 = @(9, 9)
    32 = @(7, 7)
    33 = @(7, 7)
    34 = @(6, 6)
    35 = @(5, 5)
    36 = @(6, 6)
    37 = @(7, 7)
    38 = @(7, 7)
    39 = @(5, 5)
    40 = @(3, 3)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 9).Value = $vals[0]
    $ws.Cells.Item($row, 10).Value = $vals[1]
}
